# Update column G ("K") values on Sheet1 with regenerated data.
# This reflects the commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 2
    3  = 0
    4  = 6
    5  = 2
    6  = 5
    7  = 3
    8  = 4
    9  = 3
    10 = 5
    11 = 4
    12 = 8
    13 = 3
    14 = 6
    15 = 5
    16 = 3
    17 = 4
    18 = 3
    19 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
